$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Jun_22" header that currently lives in B1 to "Jun_26" (the column is
# being kept as the "previous ratings" column and gets shifted out to column C).
$ws.Range("B1").Value = "Jun_26"

# Duplicate the whole "ratings" column (B) into a new column C, carrying values +
# styles (so the highlighted downgrade note in B11 is preserved in C11).
$ws.Range("B1:B27").Copy($ws.Range("C1:C27"))

# Column B now becomes the new "today's ratings" column (Jun_27). Reset it to the
# default "UN" (unchanged) rating for every existing firm...
$ws.Range("B1").Value = "Jun_27"
$ws.Range("B2:B27").Value = "UN"

# ...and clear the highlighted style that used to mark the downgrade note, since
# that note now lives in C11 only.
$ws.Range("B11").Style = "Normal"

# Two new analyst firms started coverage and are appended as new rows, with the
# default "UN" rating in both the new and previous ratings columns.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"

# Match the column width used for the new column C with the rest of the table.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Reflect the new selection/active cell recorded for the sheet.
$ws.Range("F12").Select()
